$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: Sampling-Medium.R
$ws.Range("A9").Value = "Sampling-Medium.R"
$ws.Range("B9").Value = 1.319
$ws.Range("C9").Value = 1.33
$ws.Range("D9").Value = 1.381
$ws.Range("E9").Formula = "=AVERAGE(B9:D9)"
$ws.Range("F9").Value = 130
$ws.Range("G9").Value = 93
$ws.Range("L9").Value = 33
$ws.Range("M9").Value = 206
$ws.Range("B9:G9").HorizontalAlignment = -4108
$ws.Range("L9:M9").HorizontalAlignment = -4108

# Row 10: Sampling-Full.R
$ws.Range("A10").Value = "Sampling-Full.R"
$ws.Range("B10").Value = 253.3
$ws.Range("C10").Value = 251.9
$ws.Range("D10").Value = 253
$ws.Range("E10").Formula = "=AVERAGE(B10:D10)"
$ws.Range("F10").Value = 4104
$ws.Range("G10").Value = 3887
$ws.Range("L10").Value = 3116
$ws.Range("M10").Value = 61042
$ws.Range("B10:G10").HorizontalAlignment = -4108
$ws.Range("L10:M10").HorizontalAlignment = -4108

# Update selection to E10
$ws.Range("E10").Select()
